$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.782.40'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.644.86'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.17'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.251'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.17'
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0842'
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('D12').Value = '1.870.36'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '1.642.59'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.17'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.526'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.62'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').Value = '26.790.65'
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  -2.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.10'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.42'
$ws.Range('E22').Value = '  +11.62%  '
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.34'
$ws.Range('E24').Value = '  -2.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.80'
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.66'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0513'
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  -3.49%  '
$ws.Range('E33').Value = '  -3.13%  '
$ws.Range('D34').Value = '1.285.31'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('E36').Value = '  +1.32%  '
$ws.Range('E37').Value = '  -4.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.538'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.34'
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('D44').Value = '1.795.91'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.26'
$ws.Range('E45').Value = '  -3.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.98'
$ws.Range('E46').Value = '  -1.92%  '
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0520'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.66'
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('E51').Value = '  -0.61%  '
